$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1) Add new slide 11: "NAT (Network Address Translation)"
# ---------------------------------------------------------------
$s11 = $p.Slides.Add(11, 2)
$s11.Shapes.Item(1).TextFrame.TextRange.Text = "NAT (Network Address Translation)`t"
$s11.Shapes.Item(2).TextFrame.TextRange.Text = "To access the Internet, one public IP address is needed, but we can use a private IP address in our private network. The idea of NAT is to allow multiple devices to access the Internet through a single public address. To achieve this, the translation of a private IP address to a public IP address is required. Network Address Translation (NAT) is a process in which one or more local IP address is translated into one or more Global IP address and vice versa in order to provide Internet access to the local hosts."

# ---------------------------------------------------------------
# 2) Add new slide 12: blank title+content
# ---------------------------------------------------------------
$s12 = $p.Slides.Add(12, 2)

# ---------------------------------------------------------------
# 3) Add a new textbox on slide 10 ("TCP/IP vs OSI") with "Internet Protocol"
# ---------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$tb = $s10.Shapes.AddTextbox(1, 318.5454330708661, 342.5454330708661, 145.45456692913385, 29.081259842519685)
$tb.Name = "Metin kutusu 5"
$tb.TextFrame.TextRange.Text = "Internet Protocol"
